$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 changes from numeric 13 to text "ITEM"
$ws.Range("A1").Value = "ITEM"

# Apply a thin box border around A1:B2 (left/right/top/bottom thin, all cells)
$rng = $ws.Range("A1:B2")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# Move the active selection to D4
$ws.Range("D4").Select()
